# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go" count refreshes (column F, and a few column G
# "lowest price" corrections) across the four sheets, plus one event
# roster change on the "全部类型" sheet: a new fhana live show is inserted
# at row 32, every entry from the old row 32 through row 37 shifts down by
# one, and the old row-37 entry (wio夏时之鸢代号鸢Only) drops off the list.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 2462
$ws1.Cells.Item(3, 6).Value = 740
$ws1.Cells.Item(4, 6).Value = 245
$ws1.Cells.Item(6, 6).Value = 698
$ws1.Cells.Item(8, 6).Value = 901
$ws1.Cells.Item(9, 6).Value = 572
$ws1.Cells.Item(10, 6).Value = 940
$ws1.Cells.Item(12, 6).Value = 129
$ws1.Cells.Item(13, 6).Value = 441
$ws1.Cells.Item(14, 6).Value = 65
$ws1.Cells.Item(16, 6).Value = 1092
$ws1.Cells.Item(17, 6).Value = 24221
$ws1.Cells.Item(17, 7).Value = "已售罄"
$ws1.Cells.Item(18, 6).Value = 2280
$ws1.Cells.Item(19, 6).Value = 144
$ws1.Cells.Item(19, 7).Value = 68
$ws1.Cells.Item(20, 6).Value = 359
$ws1.Cells.Item(22, 6).Value = 65
$ws1.Cells.Item(23, 6).Value = 354
$ws1.Cells.Item(26, 6).Value = 235
$ws1.Cells.Item(28, 6).Value = 66
$ws1.Cells.Item(29, 6).Value = 47
$ws1.Cells.Item(30, 6).Value = 354
$ws1.Cells.Item(32, 6).Value = 442

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(7, 6).Value = 265
$ws2.Cells.Item(8, 6).Value = 59
$ws2.Cells.Item(11, 6).Value = 3621
$ws2.Cells.Item(13, 6).Value = 153
$ws2.Cells.Item(15, 6).Value = 16
$ws2.Cells.Item(18, 6).Value = 21
$ws2.Cells.Item(21, 6).Value = 4126

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Cells.Item(3, 6).Value = 167
$ws3.Cells.Item(4, 6).Value = 779
$ws3.Cells.Item(5, 6).Value = 246

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(3, 6).Value = 167
$ws4.Cells.Item(4, 6).Value = 2462
$ws4.Cells.Item(5, 6).Value = 779
$ws4.Cells.Item(6, 6).Value = 740
$ws4.Cells.Item(7, 6).Value = 245
$ws4.Cells.Item(9, 6).Value = 698
$ws4.Cells.Item(14, 6).Value = 265
$ws4.Cells.Item(15, 6).Value = 246
$ws4.Cells.Item(16, 6).Value = 901
$ws4.Cells.Item(17, 6).Value = 572
$ws4.Cells.Item(18, 6).Value = 940
$ws4.Cells.Item(19, 6).Value = 129
$ws4.Cells.Item(20, 6).Value = 441
$ws4.Cells.Item(21, 6).Value = 65
$ws4.Cells.Item(23, 6).Value = 1092
$ws4.Cells.Item(24, 6).Value = 24221
$ws4.Cells.Item(24, 7).Value = "已售罄"
$ws4.Cells.Item(28, 6).Value = 153
$ws4.Cells.Item(30, 6).Value = 2280
$ws4.Cells.Item(31, 6).Value = 144
$ws4.Cells.Item(31, 7).Value = 68

# A brand-new event (fhana ONE MAN LIVE) is inserted at row 32, which
# pushes the previously-listed rows 32-37 down by one; the former row 37
# (wio夏时之鸢代号鸢Only) falls off the bottom of this shifted block, since
# row 38 onward keeps its own original identity (only its F value changes,
# handled below).

$ws4.Cells.Item(32, 3).Value = "广州·fhana ONE MAN LIVE 巡回演唱会 2024"
$ws4.Cells.Item(32, 4).Value = "流花路117号流花展贸中心5号馆 广州大麦66live house"
$ws4.Cells.Item(32, 5).Value = "2024.07.27 19:00-07.27 21:30"
$ws4.Cells.Item(32, 6).Value = 16
$ws4.Cells.Item(32, 7).Value = 330
$ws4.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87638"
$ws4.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/P9TXBIjT1718746868925.jpeg"

$ws4.Cells.Item(33, 2).NumberFormat = "@"
$ws4.Cells.Item(33, 2).Value = "2024-07-27"
$ws4.Cells.Item(33, 3).Value = "广州·原神x星穹x崩only"
$ws4.Cells.Item(33, 4).Value = "鸿盛二路巨大创意产业园 巨大产业园·智汇港"
$ws4.Cells.Item(33, 5).Value = "2024.07.27 10:00-07.27 17:00"
$ws4.Cells.Item(33, 6).Value = 359
$ws4.Cells.Item(33, 7).Value = 55
$ws4.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87184"
$ws4.Cells.Item(33, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/u67hjpFi1718160712051.jpeg"

$ws4.Cells.Item(34, 2).NumberFormat = "@"
$ws4.Cells.Item(34, 2).Value = "2024-08-02"
$ws4.Cells.Item(34, 3).Value = "广州·井草圣二 2024《夏日独白》指弹吉他音乐会"
$ws4.Cells.Item(34, 4).Value = "恩宁路265号3层 MaoLivehouse(永庆坊店)"
$ws4.Cells.Item(34, 5).Value = "2024.08.02 19:30-08.02 21:00"
$ws4.Cells.Item(34, 6).Value = 2
$ws4.Cells.Item(34, 7).Value = 260
$ws4.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86940"
$ws4.Cells.Item(34, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/iNGVydXM1717644835981.jpeg"

$ws4.Cells.Item(35, 3).Value = "广州·【暑期5折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会"
$ws4.Cells.Item(35, 4).Value = "东风中路299号 广州中山纪念堂"
$ws4.Cells.Item(35, 5).Value = "2024.08.03 20:00-08.03 21:40"
$ws4.Cells.Item(35, 6).Value = 29
$ws4.Cells.Item(35, 7).Value = 50
$ws4.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85917"
$ws4.Cells.Item(35, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/won43hte1715675570347.jpeg"

$ws4.Cells.Item(36, 3).Value = "广州·漫潮动漫游戏嘉年华"
$ws4.Cells.Item(36, 4).Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Cells.Item(36, 5).Value = "2024.08.03 10:00-08.04 17:00"
$ws4.Cells.Item(36, 6).Value = 354
$ws4.Cells.Item(36, 7).Value = 60
$ws4.Cells.Item(36, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86483"
$ws4.Cells.Item(36, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/zd0ENyNu1716775206950.jpeg"

$ws4.Cells.Item(37, 2).NumberFormat = "@"
$ws4.Cells.Item(37, 2).Value = "2024-08-03"
$ws4.Cells.Item(37, 3).Value = "广州·马娘only2024part2"
$ws4.Cells.Item(37, 4).Value = "芳村大道下市直街1号信义会馆21栋(近白鹅潭风情酒吧街) 信义会馆-21栋"
$ws4.Cells.Item(37, 5).Value = "2024.08.03 10:00-08.03 19:00"
$ws4.Cells.Item(37, 6).Value = 209
$ws4.Cells.Item(37, 7).Value = 79.90000000000001
$ws4.Cells.Item(37, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87036"
$ws4.Cells.Item(37, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/SGOLSBKb1717581022213.png"

# Rows 38+ keep their own identity; only the "want-to-go" counts move.
$ws4.Cells.Item(38, 6).Value = 235
$ws4.Cells.Item(40, 6).Value = 21
$ws4.Cells.Item(41, 6).Value = 66
$ws4.Cells.Item(42, 6).Value = 47
$ws4.Cells.Item(46, 6).Value = 442
$ws4.Cells.Item(48, 6).Value = 4126
